$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.0275626423690205
$ws.Range("J2").Value = 0.02947285804177352
$ws.Range("M2").Value = 45.1830845
$ws.Range("N2").Value = 90.366169
$ws.Range("O2").Value = 0.2982772948921854
$ws.Range("P2").Value = 0.2359735829156887
$ws.Range("Q2").Value = 0.06560583869400001
$ws.Range("R2").Value = 0.393635032164
$ws.Range("S2").Value = 0.008221310405912173
$ws.Range("T2").Value = 0.006954815910882766
$ws.Range("I3").Value = 0.0275626423690205
$ws.Range("J3").Value = 0.02947285804177352
$ws.Range("N3").Value = 73.46982600000001
$ws.Range("O3").Value = 0.1616709822417395
$ws.Range("P3").Value = 0.1918520865636367
$ws.Range("Q3").Value = 0.03555939578400001
$ws.Range("R3").Value = 0.3200345620560001
$ws.Range("S3").Value = 0.004456079464977332
$ws.Range("T3").Value = 0.00565442931230811
$ws.Range("I4").Value = 0.0275626423690205
$ws.Range("J4").Value = 0.02947285804177352
$ws.Range("M4").Value = 20.755341
$ws.Range("N4").Value = 62.26602299999999
$ws.Range("O4").Value = 0.1370169176485697
$ws.Range("P4").Value = 0.1625955454769879
$ws.Range("Q4").Value = 0.030136755132
$ws.Range("R4").Value = 0.271230796188
$ws.Range("S4").Value = 0.00377654829965306
$ws.Range("T4").Value = 0.004792155430067995
$ws.Range("I5").Value = 0.0275626423690205
$ws.Range("J5").Value = 0.02947285804177352
$ws.Range("M5").Value = 26.3069545
$ws.Range("N5").Value = 52.613909
$ws.Range("O5").Value = 0.173666037012409
$ws.Range("P5").Value = 0.1373909368441856
$ws.Range("Q5").Value = 0.03819769793400001
$ws.Range("R5").Value = 0.229186187604
$ws.Range("S5").Value = 0.004786694869818107
$ws.Range("T5").Value = 0.004049303577834952
$ws.Range("I6").Value = 0.0275626423690205
$ws.Range("J6").Value = 0.02947285804177352
$ws.Range("M6").Value = 17.34473466666667
$ws.Range("N6").Value = 52.034204
$ws.Range("O6").Value = 0.1145017121838161
$ws.Range("P6").Value = 0.1358771505744131
$ws.Range("Q6").Value = 0.02518455473600001
$ws.Range("R6").Value = 0.226660992624
$ws.Range("S6").Value = 0.003155969743563042
$ws.Range("T6").Value = 0.004004687970000362
$ws.Range("I7").Value = 0.0275626423690205
$ws.Range("J7").Value = 0.02947285804177352
$ws.Range("M7").Value = 17.400077
$ws.Range("N7").Value = 52.200231
$ws.Range("O7").Value = 0.1148670560212801
$ws.Range("P7").Value = 0.136310697625088
$ws.Range("Q7").Value = 0.025264911804
$ws.Range("R7").Value = 0.227384206236
$ws.Range("S7").Value = 0.003166039585096786
$ws.Range("T7").Value = 0.004017465840679334
$ws.Range("I8").Value = 0.1944381169324222
$ws.Range("J8").Value = 0.1386090380724913
$ws.Range("M8").Value = 45.1830845
$ws.Range("N8").Value = 90.366169
$ws.Range("O8").Value = 0.2982772948921854
$ws.Range("P8").Value = 0.2359735829156887
$ws.Range("Q8").Value = 0.4628103345335
$ws.Range("R8").Value = 1.851241338134
$ws.Range("S8").Value = 0.05799647554253332
$ws.Range("T8").Value = 0.03270807133846288
$ws.Range("I9").Value = 0.1944381169324222
$ws.Range("J9").Value = 0.1386090380724913
$ws.Range("N9").Value = 73.46982600000001
$ws.Range("O9").Value = 0.1616709822417395
$ws.Range("P9").Value = 0.1918520865636367
$ws.Range("S9").Value = 0.0314350013496989
$ws.Range("T9").Value = 0.02659243317078602
$ws.Range("I10").Value = 0.1944381169324222
$ws.Range("J10").Value = 0.1386090380724913
$ws.Range("M10").Value = 20.755341
$ws.Range("N10").Value = 62.26602299999999
$ws.Range("O10").Value = 0.1370169176485697
$ws.Range("P10").Value = 0.1625955454769879
$ws.Range("Q10").Value = 0.212596957863
$ws.Range("R10").Value = 1.275581747178
$ws.Range("S10").Value = 0.02664131145547265
$ws.Range("T10").Value = 0.02253721215343731
$ws.Range("I11").Value = 0.1944381169324222
$ws.Range("J11").Value = 0.1386090380724913
$ws.Range("M11").Value = 26.3069545
$ws.Range("N11").Value = 52.613909
$ws.Range("O11").Value = 0.173666037012409
$ws.Range("P11").Value = 0.1373909368441856
$ws.Range("Q11").Value = 0.2694621349435
$ws.Range("R11").Value = 1.077848539774
$ws.Range("S11").Value = 0.03376729721180914
$ws.Range("T11").Value = 0.01904362559585097
$ws.Range("I12").Value = 0.1944381169324222
$ws.Range("J12").Value = 0.1386090380724913
$ws.Range("M12").Value = 17.34473466666667
$ws.Range("N12").Value = 52.034204
$ws.Range("O12").Value = 0.1145017121838161
$ws.Range("P12").Value = 0.1358771505744131
$ws.Range("Q12").Value = 0.1776621171906667
$ws.Range("R12").Value = 1.065972703144
$ws.Range("S12").Value = 0.02226349730255939
$ws.Range("T12").Value = 0.01883380113715046
$ws.Range("I13").Value = 0.1944381169324222
$ws.Range("J13").Value = 0.1386090380724913
$ws.Range("M13").Value = 17.400077
$ws.Range("N13").Value = 52.200231
$ws.Range("O13").Value = 0.1148670560212801
$ws.Range("P13").Value = 0.136310697625088
$ws.Range("Q13").Value = 0.178228988711
$ws.Range("R13").Value = 1.069373932266
$ws.Range("S13").Value = 0.02233453407034875
$ws.Range("T13").Value = 0.01889389467680368
$ws.Range("G14").Value = 0.040985
$ws.Range("H14").Value = 0.122955
$ws.Range("I14").Value = 0.7779992406985573
$ws.Range("J14").Value = 0.8319181038857351
$ws.Range("M14").Value = 45.1830845
$ws.Range("N14").Value = 90.366169
$ws.Range("O14").Value = 0.2982772948921854
$ws.Range("P14").Value = 0.2359735829156887
$ws.Range("Q14").Value = 1.8518287182325
$ws.Range("R14").Value = 11.110972309395
$ws.Range("S14").Value = 0.2320595089437399
$ws.Range("T14").Value = 0.196310695666343
$ws.Range("G15").Value = 0.040985
$ws.Range("H15").Value = 0.122955
$ws.Range("I15").Value = 0.7779992406985573
$ws.Range("J15").Value = 0.8319181038857351
$ws.Range("N15").Value = 73.46982600000001
$ws.Range("O15").Value = 0.1616709822417395
$ws.Range("P15").Value = 0.1918520865636367
$ws.Range("Q15").Value = 1.00372027287
$ws.Range("R15").Value = 9.033482455830001
$ws.Range("S15").Value = 0.1257799014270633
$ws.Range("T15").Value = 0.1596052240805426
$ws.Range("G16").Value = 0.040985
$ws.Range("H16").Value = 0.122955
$ws.Range("I16").Value = 0.7779992406985573
$ws.Range("J16").Value = 0.8319181038857351
$ws.Range("M16").Value = 20.755341
$ws.Range("N16").Value = 62.26602299999999
$ws.Range("O16").Value = 0.1370169176485697
$ws.Range("P16").Value = 0.1625955454769879
$ws.Range("Q16").Value = 0.8506576508849999
$ws.Range("R16").Value = 7.655918857964998
$ws.Range("S16").Value = 0.106599057893444
$ws.Range("T16").Value = 0.1352661778934826
$ws.Range("G17").Value = 0.040985
$ws.Range("H17").Value = 0.122955
$ws.Range("I17").Value = 0.7779992406985573
$ws.Range("J17").Value = 0.8319181038857351
$ws.Range("M17").Value = 26.3069545
$ws.Range("N17").Value = 52.613909
$ws.Range("O17").Value = 0.173666037012409
$ws.Range("P17").Value = 0.1373909368441856
$ws.Range("Q17").Value = 1.0781905301825
$ws.Range("R17").Value = 6.469143181094999
$ws.Range("S17").Value = 0.1351120449307818
$ws.Range("T17").Value = 0.1142980076704996
$ws.Range("G18").Value = 0.040985
$ws.Range("H18").Value = 0.122955
$ws.Range("I18").Value = 0.7779992406985573
$ws.Range("J18").Value = 0.8319181038857351
$ws.Range("M18").Value = 17.34473466666667
$ws.Range("N18").Value = 52.034204
$ws.Range("O18").Value = 0.1145017121838161
$ws.Range("P18").Value = 0.1358771505744131
$ws.Range("Q18").Value = 0.7108739503133333
$ws.Range("R18").Value = 6.39786555282
$ws.Range("S18").Value = 0.08908224513769371
$ws.Range("T18").Value = 0.1130386614672623
$ws.Range("G19").Value = 0.040985
$ws.Range("H19").Value = 0.122955
$ws.Range("I19").Value = 0.7779992406985573
$ws.Range("J19").Value = 0.8319181038857351
$ws.Range("M19").Value = 17.400077
$ws.Range("N19").Value = 52.200231
$ws.Range("O19").Value = 0.1148670560212801
$ws.Range("P19").Value = 0.136310697625088
$ws.Range("Q19").Value = 0.713142155845
$ws.Range("R19").Value = 6.418279402604999
$ws.Range("S19").Value = 0.08936648236583455
$ws.Range("T19").Value = 0.113399337107605
